$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.531.86'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.10%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.727.77'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.59%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9992'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.73'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +2.22%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.08%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4814'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +2.27%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2665'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.71%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06190'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.06%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.732.36'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.84%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07176'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.67%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.58'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.31%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6099'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +2.57%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.40%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9996'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.09%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.531.25'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9994'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.13%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000006933'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.03%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.51'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.11%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.955.01'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.96%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.523'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.803'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.250'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.35%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '137.01'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.60%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.34'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.16%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.778'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.02%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.402'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.21%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '107.21'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.54%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.46%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08032'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +3.67%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.685'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.23%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.63%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.13%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9987'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +2.57%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6267'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.57%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9115'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.89%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.072'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +7.88%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.56%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.003'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.23%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '102.64'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -10.53%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.72%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.587'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3873'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.45%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.968'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +10.97%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1183'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.71%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05361'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.82%  '
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.821'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.22%  '
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '30.51'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.62%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.256'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +3.32%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '51.21'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.11%  '
